# Celerio v3.0.101 regen: "document.accountId" -> "document.owner" (with printer),
# and a new "owner" search-criteria row appended to the Search sheet.

$wb = $excel.ActiveWorkbook

# --- "List" sheet: rename the accountId column to owner -----------------
$list = $wb.Worksheets.Item("List")
$list.Range("B1").Value = '${msg.getProperty(''document_owner'')}'
$list.Range("B2").Value = '${printer.print(document.owner)}'

# --- "Search" sheet: add a new row for the owner search criterion -------
$search = $wb.Worksheets.Item("Search")
$search.Cells.Item(5, 1).Value = '${msg.getProperty(''document_owner'')}'
$search.Cells.Item(5, 2).Value = '${owner}'
